$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "273.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.98%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.07%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.879"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.22%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06311"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.54%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.898"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.47%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.353"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.57%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.237"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "29.47%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8731"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.46%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1458"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05150"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.27%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07360"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.11%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03045"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09045"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001595"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.04%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006319"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.63%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006033"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.96%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.454"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.272"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.15%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1325"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.12%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.929"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.78%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04400"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.90%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001176"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.04%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004402"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.94%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001694"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.54%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04030"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.31%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006701"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.00%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.00%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002100"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.46%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01175"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-16.78%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005313"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.92%"
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.02000"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-33.03%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.665"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "15.68%"
